# Weekly fruit/vegetable update: insert 3 new price rows (new week) above the
# existing "Femacal de La Calera - Kiwi / Hayward" block that starts at row 630,
# pushing the existing rows 630-664 down to 633-667.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 630..632 (shifts existing data down, copies the row
# formatting - e.g. the date style on column D - from the row above, same as
# Excel's native Insert behaviour).
$ws.Range("A630:T632").Insert()

# Row 630: Hayward / Especial
$ws.Range("A630").Value = 3
$ws.Range("B630").Value = "Femacal de La Calera"
$ws.Range("C630").Value = "Coquimbo"
$ws.Range("D630").Value = 44746
$ws.Range("E630").Value = 5
$ws.Range("F630").Value = "Fruta"
$ws.Range("G630").Value = 100101
$ws.Range("H630").Value = "Berries"
$ws.Range("I630").Value = 100101007
$ws.Range("J630").Value = "Kiwi"
$ws.Range("K630").Value = "Hayward"
$ws.Range("L630").Value = "Especial"
$ws.Range("M630").Value = 70
$ws.Range("N630").Value = 8000
$ws.Range("O630").Value = 8000
$ws.Range("P630").Value = 8000
$ws.Range("Q630").Value = "`$/bandeja 10 kilos"
$ws.Range("R630").Value = "Región de O'Higgins"
$ws.Range("S630").Value = 800
$ws.Range("T630").Value = 10

# Row 631: Hayward / Primera
$ws.Range("A631").Value = 3
$ws.Range("B631").Value = "Femacal de La Calera"
$ws.Range("C631").Value = "Coquimbo"
$ws.Range("D631").Value = 44746
$ws.Range("E631").Value = 5
$ws.Range("F631").Value = "Fruta"
$ws.Range("G631").Value = 100101
$ws.Range("H631").Value = "Berries"
$ws.Range("I631").Value = 100101007
$ws.Range("J631").Value = "Kiwi"
$ws.Range("K631").Value = "Hayward"
$ws.Range("L631").Value = "Primera"
$ws.Range("M631").Value = 85
$ws.Range("N631").Value = 7000
$ws.Range("O631").Value = 7000
$ws.Range("P631").Value = 7000
$ws.Range("Q631").Value = "`$/bandeja 10 kilos"
$ws.Range("R631").Value = "Región de O'Higgins"
$ws.Range("S631").Value = 700
$ws.Range("T631").Value = 10

# Row 632: Hayward / Segunda
$ws.Range("A632").Value = 3
$ws.Range("B632").Value = "Femacal de La Calera"
$ws.Range("C632").Value = "Coquimbo"
$ws.Range("D632").Value = 44746
$ws.Range("E632").Value = 5
$ws.Range("F632").Value = "Fruta"
$ws.Range("G632").Value = 100101
$ws.Range("H632").Value = "Berries"
$ws.Range("I632").Value = 100101007
$ws.Range("J632").Value = "Kiwi"
$ws.Range("K632").Value = "Hayward"
$ws.Range("L632").Value = "Segunda"
$ws.Range("M632").Value = 87
$ws.Range("N632").Value = 6000
$ws.Range("O632").Value = 6000
$ws.Range("P632").Value = 6000
$ws.Range("Q632").Value = "`$/bandeja 10 kilos"
$ws.Range("R632").Value = "Región de O'Higgins"
$ws.Range("S632").Value = 600
$ws.Range("T632").Value = 10
